# Update the NFL current-week schedule table (Sheet1, A1:C16) with the
# new week's matchups, times, and stadium locations. Adds one extra game
# row (row 16) versus the previous week's 15-game table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out last week's games first so the table is rebuilt fresh for
# the new week (the new slate has one more game than the old one).
$ws.Range("A2:C15").ClearContents()

$data = @(
    @("Denver  @  New Orleans",            "8:15 PM", "Caesars Superdome"),
    @("New England  vs.  Jacksonville",    "9:30 AM", "Wembley Stadium"),
    @("Philadelphia  @  NY Giants",        "1:00 PM", "MetLife Stadium"),
    @("Detroit  @  Minnesota",             "1:00 PM", "U.S. Bank Stadium"),
    @("Miami  @  Indianapolis",            "1:00 PM", "Lucas Oil Stadium"),
    @("Houston  @  Green Bay",             "1:00 PM", "Lambeau Field"),
    @("Cincinnati  @  Cleveland",          "1:00 PM", "Huntington Bank Field"),
    @("Tennessee  @  Buffalo",             "1:00 PM", "Highmark Stadium"),
    @("Seattle  @  Atlanta",               "1:00 PM", "Mercedes-Benz Stadium"),
    @("Las Vegas  @  LA Rams",             "4:05 PM", "SoFi Stadium"),
    @("Carolina  @  Washington",           "4:05 PM", "Northwest Stadium"),
    @("Kansas City  @  San Francisco",     "4:25 PM", "Levi's Stadium"),
    @("NY Jets  @  Pittsburgh",            "8:20 PM", "Acrisure Stadium"),
    @("Baltimore  @  Tampa Bay",           "8:15 PM", "Raymond James Stadium"),
    @("LA Chargers  @  Arizona",           "9:00 PM", "State Farm Stadium")
)

# Write column-by-column (all Teams, then all Time, then all Location) so
# the shared-string table is populated in the same grouped order as the
# source workbook.
for ($col = 1; $col -le 3; $col++) {
    for ($i = 0; $i -lt $data.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, $col).Value = $data[$i][$col - 1]
    }
}
